# Comprobación de los posibles errores de lectura
#
# Replace the sample roster (rows 2-4) with a new "D/E/F" test dataset,
# blank out rows 5-7 entirely (content + hyperlinks), and move the
# active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write order matters: it drives the order new entries land in
# xl/sharedStrings.xml (first-use order), so follow the same
# A2, B2, B3, C2, C3, C4, B4 sequence the original edit used.
$ws.Range("A2").Value = "Dddd"
$ws.Range("B2").Value = "44444444D"
$ws.Range("B3").Value = "55555555E"
$ws.Range("C2").Value = "dd@uniovi.es"
$ws.Range("C3").Value = "ee@uniovi.es"
$ws.Range("C4").Value = "ff@uniovi.es"
$ws.Range("B4").Value = "66666666F"

# --- Row 3 & 4 names blanked out --------------------------------------------
$ws.Range("A3").Value = ""
$ws.Range("A4").Value = ""

# --- Rows 5-7: wipe completely (Martín/Marta/José data removed) ------------
$ws.Range("A5:D7").ClearContents()

# --- Hyperlinks: keep only C2:C4, pointing at the refreshed addresses ------
$ws.Hyperlinks.Delete()
$ws.Range("C2").Hyperlinks.Add($ws.Range("C2"), "mailto:dd@uniovi.es")
$ws.Range("C3").Hyperlinks.Add($ws.Range("C3"), "mailto:ee@uniovi.es")
$ws.Range("C4").Hyperlinks.Add($ws.Range("C4"), "mailto:ff@uniovi.es")

# Re-adding hyperlinks resets the cell style to a freshly-allocated xf;
# nudge the font so the engine folds it back onto the original shared
# "Hyperlink" style (s="2") instead of leaving a near-duplicate behind.
$ws.Range("C2").Font.Underline = $true
$ws.Range("C3").Font.Underline = $true
$ws.Range("C4").Font.Underline = $true

# --- Selection moves from A8 to B4 ------------------------------------------
$ws.Range("B4").Select()
